$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.262.65"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +1.13%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.651.16"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +0.11%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.44%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "217.59"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -0.07%  "
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +2.12%  "
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -0.47%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.256"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +1.05%  "
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +1.23%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "20.00"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +1.07%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0848"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +0.28%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.882.69"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +0.10%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.657.73"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +0.54%  "
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -0.02%  "
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +2.61%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "67.72"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +1.55%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "27.252.47"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +1.01%  "
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +1.04%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "220.04"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +0.02%  "
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -0.45%  "
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +3.01%  "
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +5.67%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.43"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +0.55%  "
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +0.51%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "146.89"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +0.40%  "
$ws.Range("B26").NumberFormat = "@"
$ws.Range("B26").Value = "BinanceUSD"
$ws.Range("C26").NumberFormat = "@"
$ws.Range("C26").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.00"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -0.30%  "
$ws.Range("B27").NumberFormat = "@"
$ws.Range("B27").Value = "Cosmos"
$ws.Range("C27").NumberFormat = "@"
$ws.Range("C27").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.54"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +2.09%  "
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -0.01%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.82"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -0.52%  "
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -0.35%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.18"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -0.46%  "
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +0.04%  "
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +1.77%  "
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +1.51%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.264.14"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +1.10%  "
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -0.01%  "
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +1.07%  "
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +2.83%  "
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +1.81%  "
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -0.45%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.809"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -0.16%  "
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +1.79%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.22"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +5.68%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.792.27"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -0.07%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "62.13"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +1.31%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "91.81"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +0.28%  "
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +0.44%  "
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +1.04%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0513"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -0.28%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.65"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +0.44%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0970"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -0.64%  "
